$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# A new localization key ("strChkAbsoluteIntegral") was added to the English
# source strings and needs a row in this translation sheet (Slovak
# translation left blank, to be filled in later). It sorts alphabetically
# right before the existing "strChkComputeDerivative" row (row 9), so insert
# a fresh row there and push everything else down by one.
$ws.Rows.Item(9).Insert()

$ws.Range("B9").Value = "strChkAbsoluteIntegral"
$ws.Range("C9").Value = "In ""settings"" form, tab ""Integration"""
$ws.Range("D9").Value = "Compute the absolute-value integral?"
$ws.Range("E9").Value = ""

# Grow the "Tabla13" table (and its autofilter) so the new row is included.
$tbl = $ws.ListObjects.Item(1)
$tbl.Resize($ws.Range("B2:E168"))
